$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 193 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [double]) {
        $cell.Value2 = -1 * $val
    }
}
